$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 10-13 hold four species-observation records whose per-record fields
# (id, taxon info, coordinates, comment, biotope, substrate description...)
# need to rotate down by one row: the data that was in row 13 moves to row 10,
# row 10's data moves to row 11, row 11's moves to row 12, and row 12's moves
# to row 13. Columns that are identical across all four rows (dates, location
# names, county/municipality, reporter, etc.) are left untouched.

$cols = @("A","B","E","F","G","H","M","Q","R","AC","AH","AM","AO")

function Get-RowSnapshot($row) {
    $snap = @{}
    foreach ($col in $cols) {
        $snap[$col] = $ws.Range("$col$row").Value()
    }
    return $snap
}

function Set-RowSnapshot($row, $snap) {
    foreach ($col in $cols) {
        $v = $snap[$col]
        if ($null -eq $v) {
            $ws.Range("$col$row").Value = ""
        } else {
            $ws.Range("$col$row").Value = $v
        }
    }
}

$row10 = Get-RowSnapshot 10
$row11 = Get-RowSnapshot 11
$row12 = Get-RowSnapshot 12
$row13 = Get-RowSnapshot 13

Set-RowSnapshot 10 $row13
Set-RowSnapshot 11 $row10
Set-RowSnapshot 12 $row11
Set-RowSnapshot 13 $row12
